$d = $word.ActiveDocument

# --- Change 1: strip the direct character formatting (font/bold/color/size)
# that had been applied to the empty paragraph mark at the very end of the
# document, while keeping the paragraph-mark shading (moving it so it only
# lives on the run properties, not duplicated on the paragraph properties).
$lastParaIndex = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($lastParaIndex)
$r = $p.Range
$r.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml' w14:paraId='7121AFF7' w14:textId='77777777' w:rsidR='00102B4B' w:rsidRDefault='00102B4B' w:rsidP='00473C56'><w:pPr><w:rPr><w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/></w:rPr></w:pPr><w:bookmarkStart w:id='2' w:name='_GoBack'/><w:bookmarkEnd w:id='2'/></w:p>")

# --- Change 2: the Weeks timeline table's "9" column header becomes "8"
# (Section 8 marked complete). Edit the single character directly in that
# table cell (rather than Find/Replace, which in this runtime searches from
# the top of the story instead of staying within the supplied range) so no
# other occurrence of the digit anywhere else in the document is touched.
$tables = $d.Tables
for ($i = 1; $i -le $tables.Count; $i++) {
    $tbl = $tables.Item($i)
    if ($tbl.Range.Text -match "Weeks") {
        $cell = $tbl.Cell(2, 10)
        $ch = $cell.Range.Characters.Item(1)
        if ($ch.Text -eq "9") {
            $ch.Text = "8"
        }
        break
    }
}
